# Insert a new weekly price record at row 12 (Comercializadora del Agro de
# Limarí - Haba), pushing the existing records (old rows 12..108) down by
# one row to 13..109, and extending the used range to A1:R109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12 (carries formatting from
# the row it was inserted at, e.g. the date-formatted style on column D).
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new weekly entry.
$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 45168
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112026
$ws.Cells.Item(12, 7).Value = "Haba"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 1100
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 9000
$ws.Cells.Item(12, 13).Value = 8500
$ws.Cells.Item(12, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 340
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
